$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain text values in the
# original workbook (inline strings), including numeric-looking prices such
# as "1.001". Assigning a numeric-looking string straight to .Value would be
# auto-converted to a number by Excel, so the target cells are temporarily
# formatted as Text ("@") before the values are written, then the cell style
# is reset back to Normal (no explicit style), matching the source workbook.
$targetRange = $ws.Range("D2:E51")
$targetRange.NumberFormat = "@"

$ws.Range('D2').Value = '30.598.00'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '1.922.33'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '248.08'
$ws.Range('E5').Value = '  +3.33%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '0.4735'
$ws.Range('E7').Value = '  -0.67%  '
$ws.Range('D8').Value = '0.2903'
$ws.Range('E8').Value = '  +1.13%  '
$ws.Range('D9').Value = '0.06844'
$ws.Range('E9').Value = '  +4.03%  '
$ws.Range('D10').Value = '105.47'
$ws.Range('E10').Value = '  -1.89%  '
$ws.Range('D11').Value = '18.39'
$ws.Range('E11').Value = '  -3.65%  '
$ws.Range('D12').Value = '1.921.91'
$ws.Range('D13').Value = '0.07716'
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').Value = '5.352'
$ws.Range('E14').Value = '  +4.00%  '
$ws.Range('D15').Value = '0.6694'
$ws.Range('E15').Value = '  +1.35%  '
$ws.Range('D16').Value = '290.43'
$ws.Range('E16').Value = '  -5.79%  '
$ws.Range('D17').Value = '30.604.65'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').Value = '0.000007638'
$ws.Range('E18').Value = '  +1.33%  '
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = '12.94'
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('D21').Value = '5.542'
$ws.Range('E21').Value = '  +4.33%  '
$ws.Range('D22').Value = '2.172.41'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '6.445'
$ws.Range('E25').Value = '  +3.15%  '
$ws.Range('D26').Value = '167.79'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '20.78'
$ws.Range('E27').Value = '  +2.85%  '
$ws.Range('D28').Value = '2.125'
$ws.Range('E28').Value = '  +4.28%  '
$ws.Range('D29').Value = '0.1073'
$ws.Range('E29').Value = '  -3.87%  '
$ws.Range('E30').Value = '  +3.64%  '
$ws.Range('E31').Value = '  +1.92%  '
$ws.Range('D32').Value = '4.053'
$ws.Range('E32').Value = '  +3.15%  '
$ws.Range('D33').Value = '0.05023'
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('E34').Value = '  -1.18%  '
$ws.Range('D35').Value = '1.147'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('D36').Value = '0.02068'
$ws.Range('E36').Value = '  +6.00%  '
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').Value = '2.727'
$ws.Range('E38').Value = '  -0.94%  '
$ws.Range('D39').Value = '2.686'
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('D40').Value = '111.42'
$ws.Range('E40').Value = '  +4.24%  '
$ws.Range('E41').Value = '  -0.47%  '
$ws.Range('D42').Value = '0.4442'
$ws.Range('E42').Value = '  +6.62%  '
$ws.Range('D43').Value = '0.8759'
$ws.Range('E43').Value = '  -0.32%  '
$ws.Range('D44').Value = '5.896'
$ws.Range('E44').Value = '  +1.56%  '
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').Value = '67.81'
$ws.Range('E46').Value = '  -4.10%  '
$ws.Range('D47').Value = '7.305'
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('D48').Value = '9.403'
$ws.Range('E48').Value = '  +1.47%  '
$ws.Range('D49').Value = '48.05'
$ws.Range('E49').Value = '  +13.63%  '
$ws.Range('D50').Value = '0.1244'
$ws.Range('E50').Value = '  +3.33%  '
$ws.Range('D51').Value = '34.99'
$ws.Range('E51').Value = '  +0.24%  '

$targetRange.Style = "Normal"
